# Extend the simulation results table on Sheet1 from B2:U3 out to B2:AO3.
# The new columns V..AO repeat the same "angle" sweep as B..U in row 2, and
# carry a fresh set of 0/1 "hit" results (new simulation run) in row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: angle values (V2:AO2) -----------------------------------------
$row2 = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# --- Row 3: new simulation hit results (V3:AO3) ----------------------------
$row3 = @(1,0,1,1,0,1,0,0,1,1,1,0,0,1,1,0,0,0,1,1)

# Columns V (22) through AO (41)
$startCol = 22
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}

# Update selection / view to match the new extent being inspected.
$ws.Range("AA7").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 12
$win.ScrollRow = 1
